$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45954, 7438.42928719564, 6391.59636103991, 10820, 3291.35778876279, -47.3769104248878),
    @(45955, 3137.48363510197, 4790.13534384879, 3620, 3670.20547612158, 201.680867498765),
    @(45956, 2947.53133550742, 4654.1378544482, 3620, 3569.02017186364, 191.798251096326),
    @(45957, 8574.26417960868, 7858.78058789794, 3620, 4385.46059719256, 359.343382712104),
    @(45958, 8361.2320647443, 7787.65752752391, 3620, 4185.01282495954, 348.027931353477),
    @(45959, 7224.82758554897, 7356.46952483481, 3620, 3510.48992024833, 301.956643545131),
    @(45960, 7224.82758554897, 7064.01487633581, 3620, 3510.48992024833, 289.771033191006),
    @(45961, 6107.35191671871, 5740.58271956264, 3620, 3445.87223819891, 231.935623240065),
    @(45962, 2557.16655997904, 4089.11832702722, 9004, 3796.91782072204, -46.5818271771141),
    @(45963, 2557.16655997904, 3985.9594870564, 9004, 3796.91782072204, -50.8801121758982),
    @(45964, 8571.71877057272, 7161.84250586451, 9004, 4685.5060999062, 118.472858573779),
    @(45965, 8571.71877057272, 7174.19356055169, 9004, 4685.5060999062, 118.987485852412),
    @(45966, 8571.71877057272, 7109.11682616025, 9004, 4685.5060999062, 116.275955252768),
    @(45967, 8580.99785695336, 7132.14711735837, 9004, 4687.59519381635, 117.322596298947)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
